$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("C2").Value = -0.2878265411692921
$ws.Range("D2").Value = 0.776173523335157

# Row 3
$ws.Range("C3").Value = -0.8085972365240299
$ws.Range("D3").Value = 0.4274049896869996

# Row 4
$ws.Range("C4").Value = 0.7755806466395033
$ws.Range("D4").Value = 0.446254231299974

# Row 5
$ws.Range("C5").Value = -1.70276176128932
$ws.Range("D5").Value = 0.1027015933619893

# Row 6
$ws.Range("C6").Value = -0.58368569903801
$ws.Range("D6").Value = 0.5653686652367287

# Row 7
$ws.Range("C7").Value = 1.112914149564505
$ws.Range("D7").Value = 0.2777583239318397

# Row 8
$ws.Range("C8").Value = -1.354709238399644
$ws.Range("D8").Value = 0.1892600893798149

# Row 9
$ws.Range("C9").Value = 1.179253734505126
$ws.Range("D9").Value = 0.2508934718832188

# Row 10
$ws.Range("C10").Value = -0.5699938230354963
$ws.Range("D10").Value = 0.5744572058909845

# Row 11
$ws.Range("C11").Value = -2.020947440780088
$ws.Range("D11").Value = 0.05562015887206662
$ws.Range("G11").Value = "No"
